# edit.ps1
# Applies the "data up to 11th" update to the survey-state workbook:
#  - Adds 3 new shared strings for dates "09 06 2020", "10 06 2020", "11 06 2020"
#    (written via the A-column label cells of new rows 131-133).
#  - Fills in the AR column (a state whose results arrived late) for rows 106-108.
#  - Corrects a handful of previously-published figures on rows 127-128.
#  - Populates the two already-present date rows (129, 130) with their full data,
#    and appends three brand-new date rows (131, 132, 133) - only 131 has full data,
#    132/133 are label-only placeholders (matching the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New shared-string date rows: set column A for rows 131-133 ----
$ws.Cells.Item(131,1).Value = "09 06 2020"
$ws.Cells.Item(132,1).Value = "10 06 2020"
$ws.Cells.Item(133,1).Value = "11 06 2020"

# ---- Updated values in existing cells (rows 127-128) ----
$ws.Cells.Item(127,19).Value = 0.5629552   # S127: 0.5629997 -> 0.5629552
$ws.Cells.Item(127,39).Value = 0.3600571   # AM127: 0.3600704 -> 0.3600571
$ws.Cells.Item(128,6).Value = 0.5219863   # F128: 0.5220318 -> 0.5219863
$ws.Cells.Item(128,7).Value = 0.3303089   # G128: 0.3303167 -> 0.3303089
$ws.Cells.Item(128,8).Value = 0.4171022   # H128: 0.4171816 -> 0.4171022
$ws.Cells.Item(128,12).Value = 0.3412479   # L128: 0.3412584 -> 0.3412479
$ws.Cells.Item(128,18).Value = 0.365484   # R128: 0.3655015 -> 0.365484
$ws.Cells.Item(128,19).Value = 0.5656282   # S128: 0.5656745 -> 0.5656282
$ws.Cells.Item(128,36).Value = 0.3298458   # AJ128: 0.3298704 -> 0.3298458
$ws.Cells.Item(128,39).Value = 0.3650468   # AM128: 0.3650611 -> 0.3650468
$ws.Cells.Item(128,43).Value = 0.3173244   # AQ128: 0.3173486 -> 0.3173244
$ws.Cells.Item(128,46).Value = 0.4180059   # AT128: 0.4180432 -> 0.4180059
$ws.Cells.Item(128,50).Value = 0.5030512   # AX128: 0.5031330000000001 -> 0.5030512

# ---- New cells added (AR106-108, and data rows 129-131) ----
$ws.Cells.Item(106,44).Value = 0.6060605999999999   # AR106
$ws.Cells.Item(107,44).Value = 0.6756757   # AR107
$ws.Cells.Item(108,44).Value = 0   # AR108
$ws.Cells.Item(129,2).Value = 0.4989133   # B129
$ws.Cells.Item(129,3).Value = 0.5671654   # C129
$ws.Cells.Item(129,4).Value = 0.7330367   # D129
$ws.Cells.Item(129,6).Value = 0.5199751   # F129
$ws.Cells.Item(129,7).Value = 0.3371037   # G129
$ws.Cells.Item(129,8).Value = 0.3974634   # H129
$ws.Cells.Item(129,9).Value = 0.3222541   # I129
$ws.Cells.Item(129,10).Value = 0.318717   # J129
$ws.Cells.Item(129,11).Value = 0.3594348   # K129
$ws.Cells.Item(129,12).Value = 0.3387692   # L129
$ws.Cells.Item(129,13).Value = 0.4206858   # M129
$ws.Cells.Item(129,15).Value = 0.2380167   # O129
$ws.Cells.Item(129,16).Value = 0.4602957   # P129
$ws.Cells.Item(129,17).Value = 0.296018   # Q129
$ws.Cells.Item(129,18).Value = 0.3606113   # R129
$ws.Cells.Item(129,19).Value = 0.5647226   # S129
$ws.Cells.Item(129,20).Value = 0.4635734   # T129
$ws.Cells.Item(129,21).Value = 0.3266002   # U129
$ws.Cells.Item(129,22).Value = 0.5348586   # V129
$ws.Cells.Item(129,23).Value = 0.2424493   # W129
$ws.Cells.Item(129,24).Value = 0.5614941   # X129
$ws.Cells.Item(129,25).Value = 0.21773   # Y129
$ws.Cells.Item(129,26).Value = 0.3106558   # Z129
$ws.Cells.Item(129,27).Value = 0.4295411   # AA129
$ws.Cells.Item(129,28).Value = 0.4941162   # AB129
$ws.Cells.Item(129,30).Value = 0.6071496   # AD129
$ws.Cells.Item(129,31).Value = 0.2915715   # AE129
$ws.Cells.Item(129,32).Value = 0.387626   # AF129
$ws.Cells.Item(129,33).Value = 0.7473147   # AG129
$ws.Cells.Item(129,34).Value = 0.5228258   # AH129
$ws.Cells.Item(129,35).Value = 0.2026529   # AI129
$ws.Cells.Item(129,36).Value = 0.3581723   # AJ129
$ws.Cells.Item(129,37).Value = 0.2941307   # AK129
$ws.Cells.Item(129,38).Value = 0.4141155   # AL129
$ws.Cells.Item(129,39).Value = 0.3491154   # AM129
$ws.Cells.Item(129,40).Value = 0.3454104   # AN129
$ws.Cells.Item(129,41).Value = 0.4433821   # AO129
$ws.Cells.Item(129,42).Value = 0.2645647   # AP129
$ws.Cells.Item(129,43).Value = 0.2980715   # AQ129
$ws.Cells.Item(129,45).Value = 0.3351337   # AS129
$ws.Cells.Item(129,46).Value = 0.4077499   # AT129
$ws.Cells.Item(129,47).Value = 0.4664907   # AU129
$ws.Cells.Item(129,48).Value = 0.4900366   # AV129
$ws.Cells.Item(129,49).Value = 0.3827458   # AW129
$ws.Cells.Item(129,50).Value = 0.4735745   # AX129
$ws.Cells.Item(129,51).Value = 0.4107948   # AY129
$ws.Cells.Item(129,53).Value = 0.2371759   # BA129
$ws.Cells.Item(129,54).Value = 0.2905741   # BB129
$ws.Cells.Item(129,55).Value = 0.2654873   # BC129
$ws.Cells.Item(129,56).Value = 0.3465377   # BD129
$ws.Cells.Item(129,57).Value = 0.7170943   # BE129
$ws.Cells.Item(130,2).Value = 0.4423429   # B130
$ws.Cells.Item(130,3).Value = 0.5798375   # C130
$ws.Cells.Item(130,4).Value = 0.7163267   # D130
$ws.Cells.Item(130,6).Value = 0.5110475   # F130
$ws.Cells.Item(130,7).Value = 0.3386219   # G130
$ws.Cells.Item(130,8).Value = 0.4000528   # H130
$ws.Cells.Item(130,9).Value = 0.3601287   # I130
$ws.Cells.Item(130,10).Value = 0.2383863   # J130
$ws.Cells.Item(130,11).Value = 0.3306813   # K130
$ws.Cells.Item(130,12).Value = 0.3493911   # L130
$ws.Cells.Item(130,13).Value = 0.4486603   # M130
$ws.Cells.Item(130,15).Value = 0.2645235   # O130
$ws.Cells.Item(130,16).Value = 0.5062392   # P130
$ws.Cells.Item(130,17).Value = 0.2177891   # Q130
$ws.Cells.Item(130,18).Value = 0.3442159   # R130
$ws.Cells.Item(130,19).Value = 0.5499477   # S130
$ws.Cells.Item(130,20).Value = 0.5210275   # T130
$ws.Cells.Item(130,21).Value = 0.3281358   # U130
$ws.Cells.Item(130,22).Value = 0.5821739   # V130
$ws.Cells.Item(130,23).Value = 0.2850863   # W130
$ws.Cells.Item(130,24).Value = 0.5434924   # X130
$ws.Cells.Item(130,25).Value = 0.3470019   # Y130
$ws.Cells.Item(130,26).Value = 0.3298128   # Z130
$ws.Cells.Item(130,27).Value = 0.3417861   # AA130
$ws.Cells.Item(130,28).Value = 0.4589467   # AB130
$ws.Cells.Item(130,30).Value = 0.6091556   # AD130
$ws.Cells.Item(130,31).Value = 0.2536454   # AE130
$ws.Cells.Item(130,32).Value = 0.3685354   # AF130
$ws.Cells.Item(130,33).Value = 0.727095   # AG130
$ws.Cells.Item(130,34).Value = 0.55164   # AH130
$ws.Cells.Item(130,35).Value = 0.270764   # AI130
$ws.Cells.Item(130,36).Value = 0.342259   # AJ130
$ws.Cells.Item(130,37).Value = 0.2891348   # AK130
$ws.Cells.Item(130,38).Value = 0.3891482   # AL130
$ws.Cells.Item(130,39).Value = 0.3465515   # AM130
$ws.Cells.Item(130,40).Value = 0.3453777   # AN130
$ws.Cells.Item(130,41).Value = 0.4431783   # AO130
$ws.Cells.Item(130,42).Value = 0.2342736   # AP130
$ws.Cells.Item(130,43).Value = 0.2846097   # AQ130
$ws.Cells.Item(130,45).Value = 0.3284933   # AS130
$ws.Cells.Item(130,46).Value = 0.4059799   # AT130
$ws.Cells.Item(130,47).Value = 0.4644246   # AU130
$ws.Cells.Item(130,48).Value = 0.4609496   # AV130
$ws.Cells.Item(130,49).Value = 0.3819311   # AW130
$ws.Cells.Item(130,50).Value = 0.3732232   # AX130
$ws.Cells.Item(130,51).Value = 0.3922044   # AY130
$ws.Cells.Item(130,53).Value = 0.1873934   # BA130
$ws.Cells.Item(130,54).Value = 0.2765795   # BB130
$ws.Cells.Item(130,55).Value = 0.2795722   # BC130
$ws.Cells.Item(130,56).Value = 0.3302295   # BD130
$ws.Cells.Item(130,57).Value = 0.6461973   # BE130
$ws.Cells.Item(131,2).Value = 0.5169173   # B131
$ws.Cells.Item(131,3).Value = 0.6273738   # C131
$ws.Cells.Item(131,4).Value = 0.6471934   # D131
$ws.Cells.Item(131,6).Value = 0.5468536000000001   # F131
$ws.Cells.Item(131,7).Value = 0.3592687   # G131
$ws.Cells.Item(131,8).Value = 0.3780051   # H131
$ws.Cells.Item(131,9).Value = 0.3591511   # I131
$ws.Cells.Item(131,10).Value = 0.24375   # J131
$ws.Cells.Item(131,11).Value = 0.3346359   # K131
$ws.Cells.Item(131,12).Value = 0.343297   # L131
$ws.Cells.Item(131,13).Value = 0.4721408   # M131
$ws.Cells.Item(131,15).Value = 0.2735209   # O131
$ws.Cells.Item(131,16).Value = 0.5158463   # P131
$ws.Cells.Item(131,17).Value = 0.2383149   # Q131
$ws.Cells.Item(131,18).Value = 0.3539119   # R131
$ws.Cells.Item(131,19).Value = 0.5472557   # S131
$ws.Cells.Item(131,20).Value = 0.5294053   # T131
$ws.Cells.Item(131,21).Value = 0.3517988   # U131
$ws.Cells.Item(131,22).Value = 0.6047844999999999   # V131
$ws.Cells.Item(131,23).Value = 0.2896218   # W131
$ws.Cells.Item(131,24).Value = 0.5448647   # X131
$ws.Cells.Item(131,25).Value = 0.3675081   # Y131
$ws.Cells.Item(131,26).Value = 0.3387754   # Z131
$ws.Cells.Item(131,27).Value = 0.3145624   # AA131
$ws.Cells.Item(131,28).Value = 0.4241433   # AB131
$ws.Cells.Item(131,30).Value = 0.582677   # AD131
$ws.Cells.Item(131,31).Value = 0.2363869   # AE131
$ws.Cells.Item(131,32).Value = 0.3722949   # AF131
$ws.Cells.Item(131,33).Value = 0.6126424   # AG131
$ws.Cells.Item(131,34).Value = 0.4617437   # AH131
$ws.Cells.Item(131,35).Value = 0.281419   # AI131
$ws.Cells.Item(131,36).Value = 0.3620017   # AJ131
$ws.Cells.Item(131,37).Value = 0.3186358   # AK131
$ws.Cells.Item(131,38).Value = 0.4200534   # AL131
$ws.Cells.Item(131,39).Value = 0.3577762   # AM131
$ws.Cells.Item(131,40).Value = 0.3122303   # AN131
$ws.Cells.Item(131,41).Value = 0.4756105   # AO131
$ws.Cells.Item(131,42).Value = 0.251016   # AP131
$ws.Cells.Item(131,43).Value = 0.2814056   # AQ131
$ws.Cells.Item(131,45).Value = 0.3737468   # AS131
$ws.Cells.Item(131,46).Value = 0.4677864   # AT131
$ws.Cells.Item(131,47).Value = 0.429169   # AU131
$ws.Cells.Item(131,48).Value = 0.4380661   # AV131
$ws.Cells.Item(131,49).Value = 0.4152559   # AW131
$ws.Cells.Item(131,50).Value = 0.2826557   # AX131
$ws.Cells.Item(131,51).Value = 0.3914485   # AY131
$ws.Cells.Item(131,53).Value = 0.2554315   # BA131
$ws.Cells.Item(131,54).Value = 0.3024089   # BB131
$ws.Cells.Item(131,55).Value = 0.2477588   # BC131
$ws.Cells.Item(131,56).Value = 0.3929492   # BD131
$ws.Cells.Item(131,57).Value = 0.6590323   # BE131
